# Apply updated NATMI ligand-receptor statistics (Sertad1-Ar)
# following recount of ligand/receptor-expressing cells (1 -> 3 per group).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "E2"=3; "G2"=10.13474866666667; "H2"=30.404246; "I2"=0.09605827562377289; "J2"=0.09605827562377291; "K2"=3; "M2"=0.4955973333333333; "N2"=1.486792; "O2"=0.05374044238728271; "P2"=0.05374044238728272; "Q2"=5.022754413203555; "R2"=45.204789718832; "S2"=0.00516221422698109; "T2"=0.005162214226981092;
    "E3"=3; "G3"=10.13474866666667; "H3"=30.404246; "I3"=0.09605827562377289; "J3"=0.09605827562377291; "K3"=3; "M3"=4.129416666666667; "N3"=12.38825; "O3"=0.4477761754194636; "P3"=0.4477761754194636; "Q3"=41.85060005661111; "R3"=376.6554005095; "S3"=0.04301260727620172; "T3"=0.04301260727620172;
    "E4"=3; "G4"=10.13474866666667; "H4"=30.404246; "I4"=0.09605827562377289; "J4"=0.09605827562377291; "K4"=3; "M4"=4.597041333333333; "N4"=13.791124; "O4"=0.4984833821932537; "P4"=0.4984833821932537; "Q4"=46.58985852361155; "R4"=419.308726712504; "S4"=0.04788345412059009; "T4"=0.0478834541205901;
    "E5"=3; "G5"=8.977051666666666; "H5"=26.931155; "I5"=0.08508549463310319; "J5"=0.08508549463310322; "K5"=3; "M5"=0.4955973333333333; "N5"=1.486792; "O5"=0.05374044238728271; "P5"=0.05374044238728272; "Q5"=4.449002867195555; "R5"=40.04102580476; "S5"=0.004572532122323734; "T5"=0.004572532122323737;
    "E6"=3; "G6"=8.977051666666666; "H6"=26.931155; "I6"=0.08508549463310319; "J6"=0.08508549463310322; "K6"=3; "M6"=4.129416666666667; "N6"=12.38825; "O6"=0.4477761754194636; "P6"=0.4477761754194636; "Q6"=37.06998676986111; "R6"=333.62988092875; "S6"=0.03809925737048424; "T6"=0.03809925737048426;
    "E7"=3; "G7"=8.977051666666666; "H7"=26.931155; "I7"=0.08508549463310319; "J7"=0.08508549463310322; "K7"=3; "M7"=4.597041333333333; "N7"=13.791124; "O7"=0.4984833821932537; "P7"=0.4984833821932537; "Q7"=41.26787756313556; "R7"=371.41089806822; "S7"=0.04241370514029521; "T7"=0.04241370514029523;
    "E8"=3; "G8"=16.21867033333333; "H8"=48.656011; "I8"=0.1537223621789972; "J8"=0.1537223621789972; "K8"=3; "M8"=0.4955973333333333; "N8"=1.486792; "O8"=0.05374044238728271; "P8"=0.05374044238728272; "Q8"=8.037929767412443; "R8"=72.34136790671199; "S8"=0.008261107748317403; "T8"=0.008261107748317406;
    "E9"=3; "G9"=16.21867033333333; "H9"=48.656011; "I9"=0.1537223621789972; "J9"=0.1537223621789972; "K9"=3; "M9"=4.129416666666667; "N9"=12.38825; "O9"=0.4477761754194636; "P9"=0.4477761754194636; "Q9"=66.97364758563889; "R9"=602.76282827075; "S9"=0.06883321141295695; "T9"=0.06883321141295697;
    "E10"=3; "G10"=16.21867033333333; "H10"=48.656011; "I10"=0.1537223621789972; "J10"=0.1537223621789972; "K10"=3; "M10"=4.597041333333333; "N10"=13.791124; "O10"=0.4984833821932537; "P10"=0.4984833821932537; "Q10"=74.55789789404044; "R10"=671.021081046364; "S10"=0.0766280430177228; "T10"=0.07662804301772282;
    "E11"=3; "G11"=18.97137133333333; "H11"=56.914114; "I11"=0.1798127685684043; "J11"=0.1798127685684043; "K11"=3; "M11"=0.4955973333333333; "N11"=1.486792; "O11"=0.05374044238728271; "P11"=0.05374044238728272; "Q11"=9.402161042476443; "R11"=84.619449382288; "S11"=0.00966321772974813; "T11"=0.009663217729748133;
    "E12"=3; "G12"=18.97137133333333; "H12"=56.914114; "I12"=0.1798127685684043; "J12"=0.1798127685684043; "K12"=3; "M12"=4.129416666666667; "N12"=12.38825; "O12"=0.4477761754194636; "P12"=0.4477761754194636; "Q12"=78.34069697338889; "R12"=705.0662727604999; "S12"=0.08051587380114521; "T12"=0.08051587380114522;
    "E13"=3; "G13"=18.97137133333333; "H13"=56.914114; "I13"=0.1798127685684043; "J13"=0.1798127685684043; "K13"=3; "M13"=4.597041333333333; "N13"=13.791124; "O13"=0.4984833821932537; "P13"=0.4984833821932537; "Q13"=87.21217816934845; "R13"=784.9096035241359; "S13"=0.08963367703751093; "T13"=0.08963367703751096;
    "E14"=3; "G14"=40.62112166666667; "H14"=121.863365; "I14"=0.3850115113399109; "J14"=0.3850115113399109; "K14"=3; "M14"=0.4955973333333333; "N14"=1.486792; "O14"=0.05374044238728271; "P14"=0.05374044238728272; "Q14"=20.13171957500889; "R14"=181.18547617508; "S14"=0.02069068894360312; "T14"=0.02069068894360313;
    "E15"=3; "G15"=40.62112166666667; "H15"=121.863365; "I15"=0.3850115113399109; "J15"=0.3850115113399109; "K15"=3; "M15"=4.129416666666667; "N15"=12.38825; "O15"=0.4477761754194636; "P15"=0.4477761754194636; "Q15"=167.7415368290278; "R15"=1509.67383146125; "S15"=0.1723989820402527; "T15"=0.1723989820402528;
    "E16"=3; "G16"=40.62112166666667; "H16"=121.863365; "I16"=0.3850115113399109; "J16"=0.3850115113399109; "K16"=3; "M16"=4.597041333333333; "N16"=13.791124; "O16"=0.4984833821932537; "P16"=0.4984833821932537; "Q16"=186.7369753080289; "R16"=1680.63277777226; "S16"=0.191921840356055; "T16"=0.1919218403560551;
    "E17"=3; "G17"=10.58328866666667; "H17"=31.749866; "I17"=0.1003095876558115; "J17"=0.1003095876558115; "K17"=3; "M17"=0.4955973333333333; "N17"=1.486792; "O17"=0.05374044238728271; "P17"=0.05374044238728272; "Q17"=5.245049641096889; "R17"=47.205446769872; "S17"=0.005390681616309223; "T17"=0.005390681616309225;
    "E18"=3; "G18"=10.58328866666667; "H18"=31.749866; "I18"=0.1003095876558115; "J18"=0.1003095876558115; "K18"=3; "M18"=4.129416666666667; "N18"=12.38825; "O18"=0.4477761754194636; "P18"=0.4477761754194636; "Q18"=43.70280860827778; "R18"=393.3252774745; "S18"=0.04491624351842271; "T18"=0.04491624351842272;
    "E19"=3; "G19"=10.58328866666667; "H19"=31.749866; "I19"=0.1003095876558115; "J19"=0.1003095876558115; "K19"=3; "M19"=4.597041333333333; "N19"=13.791124; "O19"=0.4984833821932537; "P19"=0.4984833821932537; "Q19"=48.65181544326489; "R19"=437.8663389893841; "S19"=0.05000266252107956; "T19"=0.05000266252107958;
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}
